$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H: "Subtype" header + "T" for every data row (3-110)
$ws.Range("H2").Value = "Subtype"
$ws.Range("H3:H110").Value = "T"

# Empty styled cell J106 (underline font) - matches where the user's
# selection ended up after formatting an empty cell
$ws.Range("J106").Font.Underline = $true

# Restore the selection to where the author left it
$ws.Range("J106").Select() | Out-Null
